# Actualización de flujos automatizados
#
# The "Cod cliente" value recorded in the visit-verification detail sheet
# (Hoja1!A2) is refreshed to reflect the latest automated run: the old
# client code 21838047 is replaced with 15979240. The cell already carries
# a text number format (style index 2 / numFmtId 49), so writing a plain
# string here keeps it stored as text, exactly like the original value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "15979240"
